# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates/inserts/deletes per the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of Pandaemonium_Profits.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2814.8408
$ws.Range("J17").Value = 2982.0244
$ws.Range("L17").Value = 8946.073199999999
$ws.Range("N17").Value = -9282.073199999999
$ws.Range("H33").Value = 209.51852
$ws.Range("I33").Value = 193.45833
$ws.Range("J33").Value = 338
$ws.Range("K33").Value = 193.45833
$ws.Range("L33").Value = 338
$ws.Range("M33").Value = 35.54167000000001
$ws.Range("N33").Value = -796
$ws.Range("H69").Value = 5120.4287
$ws.Range("J69").Value = 5120.4287
$ws.Range("L69").Value = 15361.2861
$ws.Range("N69").Value = -17109.2861
$ws.Range("H70").Value = 1846.1428
$ws.Range("J70").Value = 1990.6
$ws.Range("L70").Value = 5971.799999999999
$ws.Range("N70").Value = -6511.799999999999
$ws.Range("H72").Value = 5120.4287
$ws.Range("J72").Value = 5120.4287
$ws.Range("L72").Value = 46083.85830000001
$ws.Range("N72").Value = -54819.85830000001
$ws.Range("H73").Value = 1846.1428
$ws.Range("J73").Value = 1990.6
$ws.Range("L73").Value = 5971.799999999999
$ws.Range("N73").Value = -7843.799999999999
$ws.Range("H88").Value = 1978.2222
$ws.Range("I88").Value = 1333.3334
$ws.Range("J88").Value = 2300.6667
$ws.Range("K88").Value = 1333.3334
$ws.Range("L88").Value = 2300.6667
$ws.Range("M88").Value = -927.3334
$ws.Range("N88").Value = -3112.6667
$ws.Range("H91").Value = 1978.2222
$ws.Range("I91").Value = 1333.3334
$ws.Range("J91").Value = 2300.6667
$ws.Range("K91").Value = 1333.3334
$ws.Range("L91").Value = 2300.6667
$ws.Range("M91").Value = 70.66660000000002
$ws.Range("N91").Value = -5108.6667
$ws.Range("H101").Value = 1860.8
$ws.Range("I101").Value = 589.3333
$ws.Range("K101").Value = 1767.9999
$ws.Range("M101").Value = -145.9999
$ws.Range("H135").Value = 115386520
$ws.Range("I135").Value = 50002296
$ws.Range("J135").Value = 333333920
$ws.Range("K135").Value = 450020664
$ws.Range("L135").Value = 3000005280
$ws.Range("M135").Value = -450018129
$ws.Range("N135").Value = -3000010350
$ws.Range("H138").Value = 4967
$ws.Range("I138").Value = 1277.5667
$ws.Range("J138").Value = 7321.9575
$ws.Range("K138").Value = 3832.7001
$ws.Range("L138").Value = 21965.8725
$ws.Range("M138").Value = 1307.2999
$ws.Range("N138").Value = -32245.8725
$ws.Range("H141").Value = 1771.3914
$ws.Range("I141").Value = 1288.8206
$ws.Range("J141").Value = 4460
$ws.Range("K141").Value = 3866.4618
$ws.Range("L141").Value = 13380
$ws.Range("M141").Value = 1313.5382
$ws.Range("N141").Value = -23740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18812.725
$ws.Range("I32").Value = 20693.455
$ws.Range("K32").Value = 20693.455
$ws.Range("M32").Value = -20406.455
$ws.Range("H61").Value = 5068.5273
$ws.Range("I61").Value = 3079.3
$ws.Range("K61").Value = 3079.3
$ws.Range("M61").Value = -2867.3
$ws.Range("H74").Value = 8704.9
$ws.Range("I74").Value = 5994.9614
$ws.Range("J74").Value = 26319.5
$ws.Range("K74").Value = 5994.9614
$ws.Range("L74").Value = 26319.5
$ws.Range("M74").Value = -5120.9614
$ws.Range("N74").Value = -28067.5
$ws.Range("H77").Value = 8704.9
$ws.Range("I77").Value = 5994.9614
$ws.Range("J77").Value = 26319.5
$ws.Range("K77").Value = 29974.807
$ws.Range("L77").Value = 131597.5
$ws.Range("M77").Value = -25606.807
$ws.Range("N77").Value = -140333.5
$ws.Range("H133").Value = 32761
$ws.Range("J133").Value = 32761
$ws.Range("L133").Value = 32761
$ws.Range("N133").Value = -37821
$ws.Range("H136").Value = 5068.5273
$ws.Range("I136").Value = 3079.3
$ws.Range("K136").Value = 9237.900000000001
$ws.Range("M136").Value = -6687.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1243.2
$ws.Range("I20").Value = 1243.2
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1243.2
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -996.2
$ws.Range("N20").ClearContents()
$ws.Range("H134").Value = 3990.7
$ws.Range("I134").Value = 4791.4614
$ws.Range("J134").Value = 2503.5715
$ws.Range("K134").Value = 14374.3842
$ws.Range("L134").Value = 7510.7145
$ws.Range("M134").Value = -11839.3842
$ws.Range("N134").Value = -12580.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1700.1777
$ws.Range("I132").Value = 1313.7838
$ws.Range("J132").Value = 3487.25
$ws.Range("K132").Value = 3941.3514
$ws.Range("L132").Value = 10461.75
$ws.Range("M132").Value = -1411.3514
$ws.Range("N132").Value = -15521.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2964.5
$ws.Range("J75").Value = 3900.4
$ws.Range("L75").Value = 11701.2
$ws.Range("N75").Value = -13697.2
$ws.Range("H78").Value = 2964.5
$ws.Range("J78").Value = 3900.4
$ws.Range("L78").Value = 35103.6
$ws.Range("N78").Value = -45087.6
$ws.Range("H131").Value = 23982.477
$ws.Range("J131").Value = 51557.79
$ws.Range("L131").Value = 154673.37
$ws.Range("N131").Value = -164753.37
$ws.Range("H132").Value = 2282.7144
$ws.Range("I132").Value = 2736.8572
$ws.Range("K132").Value = 24631.7148
$ws.Range("M132").Value = -22101.7148
$ws.Range("H139").Value = 1808055
$ws.Range("I139").Value = 3523003.2
$ws.Range("J139").Value = 2846.1052
$ws.Range("K139").Value = 10569009.6
$ws.Range("L139").Value = 8538.3156
$ws.Range("M139").Value = -10563869.6
$ws.Range("N139").Value = -18818.3156

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 15000
$ws.Range("J12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -15280
$ws.Range("H64").Value = 35271
$ws.Range("J64").Value = 35271
$ws.Range("L64").Value = 35271
$ws.Range("N64").Value = -35767
$ws.Range("H67").Value = 35271
$ws.Range("J67").Value = 35271
$ws.Range("L67").Value = 35271
$ws.Range("N67").Value = -36987
$ws.Range("H80").Value = 9381.1875
$ws.Range("J80").Value = 5999.5
$ws.Range("L80").Value = 5999.5
$ws.Range("N80").Value = -7995.5
$ws.Range("H83").Value = 9381.1875
$ws.Range("J83").Value = 5999.5
$ws.Range("L83").Value = 29997.5
$ws.Range("N83").Value = -39981.5
$ws.Range("H122").Value = 8614.25
$ws.Range("I122").Value = 10881.2
$ws.Range("J122").Value = 4836
$ws.Range("K122").Value = 32643.6
$ws.Range("L122").Value = 14508
$ws.Range("M122").Value = -30193.6
$ws.Range("N122").Value = -19408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4877.5557
$ws.Range("I7").Value = 4373.1333
$ws.Range("J7").Value = 7399.6665
$ws.Range("K7").Value = 4373.1333
$ws.Range("L7").Value = 7399.6665
$ws.Range("M7").Value = -4261.1333
$ws.Range("N7").Value = -7623.6665
$ws.Range("H9").Value = 437
$ws.Range("I9").Value = 437
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 437
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -213
$ws.Range("N9").ClearContents()
$ws.Range("H82").Value = 2194
$ws.Range("J82").Value = 2897.5
$ws.Range("L82").Value = 2897.5
$ws.Range("N82").Value = -3619.5
$ws.Range("H85").Value = 2194
$ws.Range("J85").Value = 2897.5
$ws.Range("L85").Value = 2897.5
$ws.Range("N85").Value = -5393.5
$ws.Range("H126").Value = 4877.5557
$ws.Range("I126").Value = 4373.1333
$ws.Range("J126").Value = 7399.6665
$ws.Range("K126").Value = 13119.3999
$ws.Range("L126").Value = 22198.9995
$ws.Range("M126").Value = -10649.3999
$ws.Range("N126").Value = -27138.9995
$ws.Range("H136").Value = 5374.3413
$ws.Range("I136").Value = 2928.9
$ws.Range("J136").Value = 7703.3335
$ws.Range("K136").Value = 8786.700000000001
$ws.Range("L136").Value = 23110.0005
$ws.Range("M136").Value = -6236.700000000001
$ws.Range("N136").Value = -28210.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3750
$ws.Range("I81").Value = 1400
$ws.Range("J81").Value = 4085.7144
$ws.Range("K81").Value = 2800
$ws.Range("L81").Value = 8171.4288
$ws.Range("M81").Value = -1739
$ws.Range("N81").Value = -10293.4288
$ws.Range("H84").Value = 3750
$ws.Range("I84").Value = 1400
$ws.Range("J84").Value = 4085.7144
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 40857.144
$ws.Range("M84").Value = -8696
$ws.Range("N84").Value = -51465.144
$ws.Range("H96").Value = 709.1
$ws.Range("J96").Value = 748.25
$ws.Range("L96").Value = 748.25
$ws.Range("N96").Value = -3494.25
$ws.Range("H122").Value = 1929
$ws.Range("I122").Value = 1437.75
$ws.Range("J122").Value = 3501
$ws.Range("K122").Value = 4313.25
$ws.Range("L122").Value = 10503
$ws.Range("M122").Value = -1863.25
$ws.Range("N122").Value = -15403
$ws.Range("H132").Value = 1970.4166
$ws.Range("I132").Value = 822.6
$ws.Range("K132").Value = 2467.8
$ws.Range("M132").Value = 62.19999999999982
$ws.Range("H136").Value = 8904.703
$ws.Range("I136").Value = 7476.6313
$ws.Range("J136").Value = 10412.111
$ws.Range("K136").Value = 22429.8939
$ws.Range("L136").Value = 31236.333
$ws.Range("M136").Value = -19879.8939
$ws.Range("N136").Value = -36336.333
